# Update "想去人数" (F column) values on both the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 21
    $ws.Range("F4").Value = 1473
    $ws.Range("F5").Value = 16
    $ws.Range("F8").Value = 46
    $ws.Range("F9").Value = 259
}
